# Weekly fruit/vegetable price update: a new daily record is inserted
# into the "Hortaliza, Vega Modelo de Temuco - Arveja Verde" sheet.
# The new record is placed at row 15 (sheet is sorted, most-recent-first
# per market date grouping), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 15 - shifts rows 15..82 down to 16..83
# and grows the sheet dimension from A1:R82 to A1:R83 automatically.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new record's data.
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 44560
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 100112022
$ws.Range("G15").Value = "Arveja Verde"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 65
$ws.Range("K15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 15000
$ws.Range("N15").Value = "`$/saco 25 kilos"
$ws.Range("O15").Value = "Región de La Araucanía"
$ws.Range("P15").Value = 600
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
